# Fixed #366 User content is lost after two generation without edition.
#
# The user-doc zones were stored as <w:fldSimple w:instr="..."/> (a
# "simple field"). Word collapses a simple field's cached result back into
# its instr the next time the field is updated, which is how the user
# content silently disappeared on a second generation. The fix stores the
# same fields in their "complex field" form instead:
#   <w:r><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:instrText>INSTR</w:instrText></w:r>
#   <w:r><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:fldChar w:fldCharType="end"/></w:r>
#
# This script walks every field in the document and rewrites it from the
# <w:fldSimple/> shorthand into that explicit begin/instrText/separate/end
# run sequence, preserving the field's instruction text and position.

$d = $word.ActiveDocument

function Convert-ToComplexField($field) {
    $instr = $field.Code.Text.Trim()

    # Position right before the field's begin character so we can splice
    # the replacement runs in at exactly the same spot.
    $insertAt = $field.Code.Start - 1

    # Removes the <w:fldSimple/> (only the field itself, not surrounding
    # runs/paragraphs).
    $field.Delete()

    $openXml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $insertionPoint = $d.Range($insertAt, $insertAt)
    $insertionPoint.InsertXML($openXml)
}

# Walk back-to-front so deleting/inserting one field never shifts the
# index of the fields still waiting to be converted.
$fieldCount = $d.Fields.Count
for ($i = $fieldCount; $i -ge 1; $i--) {
    Convert-ToComplexField $d.Fields.Item($i)
}
